$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Sheet1")

# Clear the per-site/building breakdown counts (columns AB..AK) and the
# DIFFERENCE column (AM) for the data rows (2-20), keeping the
# PREVIOUS ACCOMPLISHMENT column (AL) intact, reflecting the most
# updated status accomplishment values as of May.
$ws.Range("AB2:AK20").ClearContents()
$ws.Range("AM2:AM20").ClearContents()
